$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the D (Price) and E (Volume) columns so that
# numeric-looking strings (e.g. "1.00", "0.388") are stored verbatim as text
# instead of being auto-converted to numbers by Excel.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range("D2").Value = "65.179.57"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "3.536.86"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "594.87"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "135.62"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("D7").Value = "3.534.15"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("D12").Value = "0.388"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "4.136.15"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").Value = "27.78"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "0.0000183"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "3.531.15"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "65.128.34"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "10.13"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "14.48"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("D22").Value = "393.16"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "3.675.62"
$ws.Range("D25").Value = "74.68"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -3.58%  "
$ws.Range("D28").Value = "7.76"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D29").Value = "1.58"
$ws.Range("E29").Value = "  +9.85%  "
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.28"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "8.41"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").Value = "3.538.11"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").Value = "24.26"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").Value = "5.31"
$ws.Range("E37").Value = "  +6.14%  "
$ws.Range("D38").Value = "7.00"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "1.58"
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("D40").Value = "168.74"
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("D41").Value = "0.0819"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "0.827"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("E43").Value = "  +4.26%  "
$ws.Range("D44").Value = "25.80"
$ws.Range("E44").Value = "  -4.00%  "
$ws.Range("D45").Value = "42.94"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").Value = "4.44"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "2.421.17"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("E51").Value = "  +5.61%  "

# Restore the original (default) cell style now that the text values are set,
# so no residual per-cell style index is left behind.
$numRange.Style = "Normal"

